$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Customer Overview")

# --- Update yearly revenue values (Column C) ---
# Force text format so the "<euro> NNN.NN" strings are kept as literal text
# instead of being auto-parsed by Excel into currency numbers.
$revRange = $ws.Range("C4:C9")
$revRange.NumberFormat = "@"
$ws.Range("C4").Value = "€ 1700.19"
$ws.Range("C5").Value = "€ 1649.48"
$ws.Range("C6").Value = "€ 1342.54"
$ws.Range("C7").Value = "€ 1460.25"
$ws.Range("C8").Value = "€ 1448.54"
$ws.Range("C9").Value = "€ 1423.82"
# Restore the default (General/Normal) cell style now that the literal text is stored
$revRange.Style = "Normal"

# --- Re-align monthly Kwh figures to correct month columns ---
# Row 2
$v_N2 = $ws.Range("N2").Value2
$ws.Range("P2").Value = $v_N2
$v_M2 = $ws.Range("M2").Value2
$ws.Range("O2").Value = $v_M2
$v_L2 = $ws.Range("L2").Value2
$ws.Range("N2").Value = $v_L2
$v_K2 = $ws.Range("K2").Value2
$ws.Range("M2").Value = $v_K2
$v_J2 = $ws.Range("J2").Value2
$ws.Range("L2").Value = $v_J2
$v_I2 = $ws.Range("I2").Value2
$ws.Range("J2").Value = $v_I2
# H2 stays the same, no move needed
$ws.Range("I2").ClearContents()
$ws.Range("K2").ClearContents()

# Row 5
$v_O5 = $ws.Range("O5").Value2
$ws.Range("P5").Value = $v_O5
$v_N5 = $ws.Range("N5").Value2
$ws.Range("O5").Value = $v_N5
$v_M5 = $ws.Range("M5").Value2
$ws.Range("N5").Value = $v_M5
$v_L5 = $ws.Range("L5").Value2
$ws.Range("M5").Value = $v_L5
$v_K5 = $ws.Range("K5").Value2
$ws.Range("L5").Value = $v_K5
$v_J5 = $ws.Range("J5").Value2
$ws.Range("K5").Value = $v_J5
$v_I5 = $ws.Range("I5").Value2
$ws.Range("J5").Value = $v_I5
$v_H5 = $ws.Range("H5").Value2
$ws.Range("I5").Value = $v_H5
# G5 stays the same, no move needed
$ws.Range("H5").ClearContents()

# Row 6
$v_O6 = $ws.Range("O6").Value2
$ws.Range("P6").Value = $v_O6
$v_N6 = $ws.Range("N6").Value2
$ws.Range("O6").Value = $v_N6
$v_M6 = $ws.Range("M6").Value2
$ws.Range("N6").Value = $v_M6
$v_L6 = $ws.Range("L6").Value2
$ws.Range("M6").Value = $v_L6
$v_K6 = $ws.Range("K6").Value2
$ws.Range("L6").Value = $v_K6
$v_J6 = $ws.Range("J6").Value2
$ws.Range("K6").Value = $v_J6
$v_I6 = $ws.Range("I6").Value2
$ws.Range("J6").Value = $v_I6
$v_H6 = $ws.Range("H6").Value2
$ws.Range("I6").Value = $v_H6
# G6 stays the same, no move needed
$ws.Range("H6").ClearContents()

# Row 7
$v_N7 = $ws.Range("N7").Value2
$ws.Range("P7").Value = $v_N7
$v_M7 = $ws.Range("M7").Value2
$ws.Range("O7").Value = $v_M7
$v_L7 = $ws.Range("L7").Value2
$ws.Range("N7").Value = $v_L7
$v_K7 = $ws.Range("K7").Value2
$ws.Range("M7").Value = $v_K7
$v_J7 = $ws.Range("J7").Value2
$ws.Range("L7").Value = $v_J7
$v_I7 = $ws.Range("I7").Value2
$ws.Range("K7").Value = $v_I7
$v_H7 = $ws.Range("H7").Value2
$ws.Range("J7").Value = $v_H7
$v_G7 = $ws.Range("G7").Value2
$ws.Range("I7").Value = $v_G7
$v_F7 = $ws.Range("F7").Value2
$ws.Range("H7").Value = $v_F7
$v_E7 = $ws.Range("E7").Value2
$ws.Range("F7").Value = $v_E7
# C7 stays the same, no move needed
$ws.Range("E7").ClearContents()
$ws.Range("G7").ClearContents()

# Row 8
$v_N8 = $ws.Range("N8").Value2
$ws.Range("P8").Value = $v_N8
$v_M8 = $ws.Range("M8").Value2
$ws.Range("O8").Value = $v_M8
$v_L8 = $ws.Range("L8").Value2
$ws.Range("N8").Value = $v_L8
$v_K8 = $ws.Range("K8").Value2
$ws.Range("M8").Value = $v_K8
$v_J8 = $ws.Range("J8").Value2
$ws.Range("K8").Value = $v_J8
$v_I8 = $ws.Range("I8").Value2
$ws.Range("J8").Value = $v_I8
$v_H8 = $ws.Range("H8").Value2
$ws.Range("I8").Value = $v_H8
$v_G8 = $ws.Range("G8").Value2
$ws.Range("H8").Value = $v_G8
# F8 stays the same, no move needed
$ws.Range("G8").ClearContents()
$ws.Range("L8").ClearContents()

# Row 9
$v_N9 = $ws.Range("N9").Value2
$ws.Range("P9").Value = $v_N9
$v_M9 = $ws.Range("M9").Value2
$ws.Range("O9").Value = $v_M9
$v_L9 = $ws.Range("L9").Value2
$ws.Range("N9").Value = $v_L9
$v_K9 = $ws.Range("K9").Value2
$ws.Range("M9").Value = $v_K9
$v_J9 = $ws.Range("J9").Value2
$ws.Range("L9").Value = $v_J9
$v_I9 = $ws.Range("I9").Value2
$ws.Range("K9").Value = $v_I9
$v_H9 = $ws.Range("H9").Value2
$ws.Range("I9").Value = $v_H9
# G9 stays the same, no move needed
$ws.Range("H9").ClearContents()
$ws.Range("J9").ClearContents()
